# Insert a new data row above the current row 342 (shifts rows 342..434
# down to 343..435, carrying their values/formatting with them) and then
# populate the newly-inserted row 342 with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(342).Insert()

$ws.Cells.Item(342, 1).Value  = 1
$ws.Cells.Item(342, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(342, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(342, 4).Value  = 45211
$ws.Cells.Item(342, 5).Value  = 15
$ws.Cells.Item(342, 6).Value  = "Fruta"
$ws.Cells.Item(342, 7).Value  = 100108
$ws.Cells.Item(342, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(342, 9).Value  = 100108006
$ws.Cells.Item(342, 10).Value = "Plátano"
$ws.Cells.Item(342, 11).Value = "Sin especificar"
$ws.Cells.Item(342, 12).Value = "Pintón"
$ws.Cells.Item(342, 13).Value = 216
$ws.Cells.Item(342, 14).Value = 25000
$ws.Cells.Item(342, 15).Value = 26000
$ws.Cells.Item(342, 16).Value = 25500
$ws.Cells.Item(342, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(342, 18).Value = "Ecuador"
$ws.Cells.Item(342, 19).Value = 1275
$ws.Cells.Item(342, 20).Value = 20
